$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)  # 展览
$ws1.Range("F3").Value = 265
$ws1.Range("F7").Value = 4875
$ws1.Range("F8").Value = 4875
$ws1.Range("F14").Value = 648
$ws1.Range("F15").Value = 4551
$ws1.Range("F18").Value = 83
$ws1.Range("F20").Value = 3580
$ws1.Range("F24").Value = 3319
$ws1.Range("F26").Value = 139
$ws1.Range("F31").Value = 188
$ws1.Range("F32").Value = 94
$ws1.Range("F37").Value = 5801
$ws1.Range("F38").Value = 909
$ws1.Range("F42").Value = 55
$ws1.Range("F43").Value = 1168
$ws1.Range("F44").Value = 535
$ws1.Range("F46").Value = 2055
$ws1.Range("F47").Value = 304

$ws2 = $wb.Worksheets.Item(2)  # 演出
$ws2.Range("F4").Value = 14
$ws2.Range("F9").Value = 43
$ws2.Range("F24").Value = 762

$ws4 = $wb.Worksheets.Item(4)  # 全部类型
$ws4.Range("F3").Value = 14
$ws4.Range("F4").Value = 265
$ws4.Range("F10").Value = 4875
$ws4.Range("F11").Value = 4875
$ws4.Range("F13").Value = 43
$ws4.Range("F17").Value = 648
$ws4.Range("F18").Value = 4551
$ws4.Range("F21").Value = 83
$ws4.Range("F23").Value = 3580
$ws4.Range("F24").Value = 3319
$ws4.Range("F26").Value = 139
$ws4.Range("F28").Value = 188
$ws4.Range("F29").Value = 94
$ws4.Range("F35").Value = 5801
$ws4.Range("F36").Value = 909
$ws4.Range("F42").Value = 55
$ws4.Range("F43").Value = 1168
$ws4.Range("F44").Value = 535
$ws4.Range("F45").Value = 2055
$ws4.Range("F46").Value = 304

